# "Better notes on slide Dropout"
#
# The Dropout slide (slide 15) has speaker notes that get expanded/
# reorganized: the four short paragraphs about random deactivation /
# p = 0.1 / epochs / model variety are merged into one paragraph (ending
# with the "Caution" sentence trimmed down), and a brand-new second
# paragraph ("Caution2: ...") is appended explaining output scaling by
# 1/(1-p).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$notes = $s.NotesPage

# Find the notes body placeholder robustly (avoid relying on a hard index).
$notesShape = $null
for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
    $candidate = $notes.Shapes.Item($i)
    if ($candidate.Name -like "*Notes Placeholder*") {
        $notesShape = $candidate
    }
}

$tr = $notesShape.TextFrame.TextRange

$para1 = "The choice of neurons to be deactivated is random. All neurons are assigned a probability p which determines their activation. A new hyperparameter!! When p = 0.1, each neuron has a 1 in 10 chance of being deactivated. At each epoch, we apply this random deactivation. That is, at each pass (forward propagation) the model will learn with a different configuration of neurons, with the neurons randomly turning on and off. This procedure effectively generates slightly different models with different neuron configurations at each iteration.Caution: Dropout is only active during model training. During tests, each neuron remains active."

$para2 = "Caution2: Consider the neurons at the output layer. During training, each neuron usually get activations only from n neurons from the hidden layer due to dropout. Now, imagine we finished the training and remove dropout. Now activations of the output neurons will be computed based on more neurons ! This is likely to put the output neurons in unusual regime, so they will produce too large absolute values, being overexcited. To avoid this, the outputs are scaled by a factor of 1/(1-p) during training."

$tr.Text = $para1 + "`n" + $para2
